# Regenerate the Handback status report timestamps (commit: "Generate Report for Handback")
$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for ec3a3e78-...md moves from
# 2016-11-09 06:29:14 -> 2016-11-09 06:31:05
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-09 06:31:05"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the
# ec3a3e78-... handoff/handback pair are refreshed
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-11-09 06:30:50"
$wsZhCn.Range("K2").Value = "2016-11-09 06:31:41"

# de-de sheet: "Correspond Handoff Datetime" (shared with Overview's G2 value) and
# "Correspond Handback DateTime" are refreshed
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-11-09 06:31:05"
$wsDeDe.Range("K2").Value = "2016-11-09 06:31:59"
